$d = $word.ActiveDocument

# Locate the (currently empty) paragraph that holds the "_GoBack" bookmark.
# It is the paragraph right after the "Mased language modeling" paragraph.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.Trim() -eq "Mased language modeling") {
        $anchorIndex = $i + 1
        break
    }
}

$bmRange = $d.Paragraphs.Item($anchorIndex).Range

function Set-MaskedIndent($rng) {
    $pf = $rng.ParagraphFormat
    $pf.LeftIndent = 66
    $pf.FirstLineIndent = 10
    $pf.CharacterUnitLeftIndent = 0
    $pf.CharacterUnitFirstLineIndent = 100
}

# Insert four new (empty) paragraphs ahead of the bookmark paragraph, in
# document order; after this the bookmark paragraph sits at
# anchorIndex + 4, and the four blanks occupy anchorIndex .. anchorIndex+3.
$bmRange.InsertParagraphBefore()
$bmRange.InsertParagraphBefore()
$bmRange.InsertParagraphBefore()
$bmRange.InsertParagraphBefore()

$p1 = $d.Paragraphs.Item($anchorIndex)
$p1.Range.InsertAfter("- Hide a selective word")
Set-MaskedIndent $p1.Range

$p2 = $d.Paragraphs.Item($anchorIndex + 1)
$p2.Range.InsertAfter("- Train model to predict the masked word")
Set-MaskedIndent $p2.Range

$p3 = $d.Paragraphs.Item($anchorIndex + 2)
$p3.Range.InsertAfter("Example:")
Set-MaskedIndent $p3.Range

$p4 = $d.Paragraphs.Item($anchorIndex + 3)
$p4.Range.InsertAfter("Original text: The quick brown fox jumps over the lazy dog")
Set-MaskedIndent $p4.Range

# The bookmark paragraph (now at anchorIndex + 4) becomes the 5th new
# paragraph: text after the bookmark is inserted first (it lands after
# bookmarkEnd because the collapsed range sits right before the
# bookmark), then the text before the bookmark is inserted (it lands
# before bookmarkStart) -- this keeps two distinct runs with the
# bookmark sitting between them, exactly like the target markup.
$p5 = $d.Paragraphs.Item($anchorIndex + 4)
$p5.Range.InsertAfter("The quick [MASK] fox jumps over the lazy dog")

$p5again = $d.Paragraphs.Item($anchorIndex + 4)
$insBefore = $d.Range($p5again.Range.Start, $p5again.Range.Start)
$insBefore.InsertBefore("Masked text: ")

$p5final = $d.Paragraphs.Item($anchorIndex + 4)
Set-MaskedIndent $p5final.Range

# Finally, add a fresh empty paragraph after the bookmark paragraph (this
# keeps the original, non-indented paragraph formatting).
$p5final.Range.InsertParagraphAfter()
